# Loan RBI, Variable Instalments
# Insert a new (blank) column before column N ("Late") on the "Repayment schedule"
# sheet, matching column M's width, then make "Repayment schedule" the active
# (selected) sheet/tab with R7 selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column at N, shifting Late/heading/Outstanding columns right.
$ws.Columns("N").Insert()

# Give the newly inserted column the same width as column M (to its left).
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Make "Repayment schedule" the active sheet/tab and select cell R7 on it.
$ws.Activate() | Out-Null
$ws.Range("R7").Select() | Out-Null
